$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new benchmark row (row 13) - GPU brute force run against zebra1
# (order matters: new shared-string entries are appended in the order cells
# are written, so write columns in the same order the source file used)
$ws.Range("B13").Value = "zebra1"
$ws.Range("C13").Value = "`$i = Get-Date`n>> .\HashSekv.exe 1 4 d85fb95cb761f5874f35ce32c305739b 6 6`n>> `$j = Get-Date`n>> `$j-`$i"
$ws.Range("D13").Value = "Brute force GPU"
$ws.Range("E13").Value = "malá a velká písmena čísla"
$ws.Range("G13").Value = "128 vláken, 32 bloků, 1000 threshold"
$ws.Range("F13").Value = "267,4 s"

# Row height matches the other benchmark rows (e.g. row 9) which use the same layout
$ws.Rows("13").RowHeight = 60

# Move viewport / selection to reflect where the author was working
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D18").Select()
